$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the weekly "Fecha" values in column D down one week for rows 180-187.
$ws.Range("D180").Value = 44448
$ws.Range("D181").Value = 44448
$ws.Range("D182").Value = 44167
$ws.Range("D183").Value = 44167
$ws.Range("D184").Value = 44238
$ws.Range("D185").Value = 44238
$ws.Range("D186").Value = 44399
$ws.Range("D187").Value = 44399

# Add two new rows (188-189) for the new week, containing the same
# Primera/Segunda quality data that used to be the latest week (dated 44400).
# Match the date number format used by the other "Fecha" cells in column D.
$ws.Range("D188:D189").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A188").Value = 11
$ws.Range("B188").Value = "Vega Monumental Concepción"
$ws.Range("C188").Value = "Bíobío"
$ws.Range("D188").Value = 44400
$ws.Range("E188").Value = 8
$ws.Range("F188").Value = 100114014
$ws.Range("G188").Value = "Betarraga"
$ws.Range("H188").Value = "Sin especificar"
$ws.Range("I188").Value = "Primera"
$ws.Range("J188").Value = 600
$ws.Range("K188").Value = 600
$ws.Range("L188").Value = 700
$ws.Range("M188").Value = 650
$ws.Range("N188").Value = "$/paquete 5 unidades"
$ws.Range("O188").Value = "Región Metropolitana"
$ws.Range("P188").Value = 130
$ws.Range("Q188").Value = 5
$ws.Range("R188").Value = "Hortaliza"

$ws.Range("A189").Value = 11
$ws.Range("B189").Value = "Vega Monumental Concepción"
$ws.Range("C189").Value = "Bíobío"
$ws.Range("D189").Value = 44400
$ws.Range("E189").Value = 8
$ws.Range("F189").Value = 100114014
$ws.Range("G189").Value = "Betarraga"
$ws.Range("H189").Value = "Sin especificar"
$ws.Range("I189").Value = "Segunda"
$ws.Range("J189").Value = 300
$ws.Range("K189").Value = 500
$ws.Range("L189").Value = 500
$ws.Range("M189").Value = 500
$ws.Range("N189").Value = "$/paquete 5 unidades"
$ws.Range("O189").Value = "Región Metropolitana"
$ws.Range("P189").Value = 100
$ws.Range("Q189").Value = 5
$ws.Range("R189").Value = "Hortaliza"
